$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before A (id); shifts old A:E -> B:F
$ws.Columns("A").Insert()

# Insert 3 new columns (_id, date, __v) between role (now E) and err (now F, will become I)
$ws.Columns("F:H").Insert()

# Header row
$ws.Range("A1").Value = "id"
$ws.Range("F1").Value = "_id"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "__v"

# Row 2
$ws.Range("A2").Value = "ov5egjxyvn6ytttmhtlt"
$ws.Range("D2").Value = "VW5OGxrS"
$ws.Range("G2").Value = 44963.47878810185
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 'E11000 duplicate key error collection: test.users index: userName_1 dup key: { userName: "test_1" }'

# Row 3
$ws.Range("A3").Value = "dqu7x73wxhtc3fhfllpl"
$ws.Range("D3").Value = "wlxva8wJ"
$ws.Range("G3").Value = 44963.478788113425
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 'E11000 duplicate key error collection: test.users index: userName_1 dup key: { userName: "test_2" }'

# Row 4
$ws.Range("A4").Value = "w8503zm8pmpa3fv0j6p0"
$ws.Range("D4").Value = "PfzwR9JA"
$ws.Range("E4").Value = "vendor"
$ws.Range("G4").Value = 44963.478788113425
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 'E11000 duplicate key error collection: test.users index: userName_1 dup key: { userName: "test_3" }'

# Apply date number format to G2, then copy format only to G3:G4 so every date cell shares one style
$ws.Range("G2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Copy()
$ws.Range("G3:G4").PasteSpecial(-4122)
